$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Efnb2"
$ws.Cells.Item(2, 3).Value = "Epha4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 20.94432133333333
$ws.Cells.Item(2, 8).Value = 62.832964
$ws.Cells.Item(2, 9).Value = 0.7396577289668299
$ws.Cells.Item(2, 10).Value = 0.7396577289668298
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 2.856403666666667
$ws.Cells.Item(2, 14).Value = 8.569211
$ws.Cells.Item(2, 15).Value = 0.235832554697756
$ws.Cells.Item(2, 16).Value = 0.235832554697756
$ws.Cells.Item(2, 17).Value = 59.82543625237822
$ws.Cells.Item(2, 18).Value = 538.428926271404
$ws.Cells.Item(2, 19).Value = 0.1744353718241879
$ws.Cells.Item(2, 20).Value = 0.1744353718241879

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Efnb2"
$ws.Cells.Item(3, 3).Value = "Epha4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 20.94432133333333
$ws.Cells.Item(3, 8).Value = 62.832964
$ws.Cells.Item(3, 9).Value = 0.7396577289668299
$ws.Cells.Item(3, 10).Value = 0.7396577289668298
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.000300666666668
$ws.Cells.Item(3, 14).Value = 21.000902
$ws.Cells.Item(3, 15).Value = 0.5779641054021444
$ws.Cells.Item(3, 16).Value = 0.5779641054021444
$ws.Cells.Item(3, 17).Value = 146.6165465926143
$ws.Cells.Item(3, 18).Value = 1319.548919333528
$ws.Cells.Item(3, 19).Value = 0.4274956176260956
$ws.Cells.Item(3, 20).Value = 0.4274956176260956

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Efnb2"
$ws.Cells.Item(4, 3).Value = "Epha4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 20.94432133333333
$ws.Cells.Item(4, 8).Value = 62.832964
$ws.Cells.Item(4, 9).Value = 0.7396577289668299
$ws.Cells.Item(4, 10).Value = 0.7396577289668298
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.255294666666666
$ws.Cells.Item(4, 14).Value = 6.765884
$ws.Cells.Item(4, 15).Value = 0.1862033399000996
$ws.Cells.Item(4, 16).Value = 0.1862033399000996
$ws.Cells.Item(4, 17).Value = 47.23561620001956
$ws.Cells.Item(4, 18).Value = 425.120545800176
$ws.Cells.Item(4, 19).Value = 0.1377267395165463
$ws.Cells.Item(4, 20).Value = 0.1377267395165463

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Efnb2"
$ws.Cells.Item(5, 3).Value = "Epha4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 2.327094666666667
$ws.Cells.Item(5, 8).Value = 6.981284
$ws.Cells.Item(5, 9).Value = 0.08218235047311259
$ws.Cells.Item(5, 10).Value = 0.08218235047311258
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 2.856403666666667
$ws.Cells.Item(5, 14).Value = 8.569211
$ws.Cells.Item(5, 15).Value = 0.235832554697756
$ws.Cells.Item(5, 16).Value = 0.235832554697756
$ws.Cells.Item(5, 17).Value = 6.647121738547112
$ws.Cells.Item(5, 18).Value = 59.824095646924
$ws.Cells.Item(5, 19).Value = 0.01938127366314048
$ws.Cells.Item(5, 20).Value = 0.01938127366314047

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Efnb2"
$ws.Cells.Item(6, 3).Value = "Epha4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.327094666666667
$ws.Cells.Item(6, 8).Value = 6.981284
$ws.Cells.Item(6, 9).Value = 0.08218235047311259
$ws.Cells.Item(6, 10).Value = 0.08218235047311258
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.000300666666668
$ws.Cells.Item(6, 14).Value = 21.000902
$ws.Cells.Item(6, 15).Value = 0.5779641054021444
$ws.Cells.Item(6, 16).Value = 0.5779641054021444
$ws.Cells.Item(6, 17).Value = 16.29036234646312
$ws.Cells.Item(6, 18).Value = 146.613261118168
$ws.Cells.Item(6, 19).Value = 0.04749844867103802
$ws.Cells.Item(6, 20).Value = 0.04749844867103801

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Efnb2"
$ws.Cells.Item(7, 3).Value = "Epha4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.327094666666667
$ws.Cells.Item(7, 8).Value = 6.981284
$ws.Cells.Item(7, 9).Value = 0.08218235047311259
$ws.Cells.Item(7, 10).Value = 0.08218235047311258
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.255294666666666
$ws.Cells.Item(7, 14).Value = 6.765884
$ws.Cells.Item(7, 15).Value = 0.1862033399000996
$ws.Cells.Item(7, 16).Value = 0.1862033399000996
$ws.Cells.Item(7, 17).Value = 5.248284190561778
$ws.Cells.Item(7, 18).Value = 47.234557715056
$ws.Cells.Item(7, 19).Value = 0.01530262813893409
$ws.Cells.Item(7, 20).Value = 0.01530262813893409

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Efnb2"
$ws.Cells.Item(8, 3).Value = "Epha4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.044818
$ws.Cells.Item(8, 8).Value = 15.134454
$ws.Cells.Item(8, 9).Value = 0.1781599205600575
$ws.Cells.Item(8, 10).Value = 0.1781599205600575
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 2.856403666666667
$ws.Cells.Item(8, 14).Value = 8.569211
$ws.Cells.Item(8, 15).Value = 0.235832554697756
$ws.Cells.Item(8, 16).Value = 0.235832554697756
$ws.Cells.Item(8, 17).Value = 14.410036632866
$ws.Cells.Item(8, 18).Value = 129.690329695794
$ws.Cells.Item(8, 19).Value = 0.04201590921042763
$ws.Cells.Item(8, 20).Value = 0.04201590921042763

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Efnb2"
$ws.Cells.Item(9, 3).Value = "Epha4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.044818
$ws.Cells.Item(9, 8).Value = 15.134454
$ws.Cells.Item(9, 9).Value = 0.1781599205600575
$ws.Cells.Item(9, 10).Value = 0.1781599205600575
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 7.000300666666668
$ws.Cells.Item(9, 14).Value = 21.000902
$ws.Cells.Item(9, 15).Value = 0.5779641054021444
$ws.Cells.Item(9, 16).Value = 0.5779641054021444
$ws.Cells.Item(9, 17).Value = 35.31524280861201
$ws.Cells.Item(9, 18).Value = 317.8371852775081
$ws.Cells.Item(9, 19).Value = 0.1029700391050108
$ws.Cells.Item(9, 20).Value = 0.1029700391050108

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Efnb2"
$ws.Cells.Item(10, 3).Value = "Epha4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.044818
$ws.Cells.Item(10, 8).Value = 15.134454
$ws.Cells.Item(10, 9).Value = 0.1781599205600575
$ws.Cells.Item(10, 10).Value = 0.1781599205600575
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.255294666666666
$ws.Cells.Item(10, 14).Value = 6.765884
$ws.Cells.Item(10, 15).Value = 0.1862033399000996
$ws.Cells.Item(10, 16).Value = 0.1862033399000996
$ws.Cells.Item(10, 17).Value = 11.377551129704
$ws.Cells.Item(10, 18).Value = 102.397960167336
$ws.Cells.Item(10, 19).Value = 0.03317397224461913
$ws.Cells.Item(10, 20).Value = 0.03317397224461913

